# "Fixed data, added availability"
# 1) Reduce the B-column prices on the VoltageRelay sheet.
# 2) Add a new "Metadata" sheet after VoltageRelay with producer/date/comment info.
# 3) Resize a handful of columns on VoltageRelay and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- 1) Updated prices -----------------------------------------------------
$ws.Range("B2").Value = 95
$ws.Range("B3").Value = 103
$ws.Range("B4").Value = 118
$ws.Range("B5").Value = 190
$ws.Range("B6").Value = 117
$ws.Range("B7").Value = 131
$ws.Range("B8").Value = 157
$ws.Range("B9").Value = 190
$ws.Range("B10").Value = 220
$ws.Range("B11").Value = 252
$ws.Range("B12").Value = 283
$ws.Range("B13").Value = 315
$ws.Range("B14").Value = 421
$ws.Range("B15").Value = 298

# ---- 2) Column width changes on VoltageRelay -------------------------------
$ws.Columns.Item(1).ColumnWidth = 34.5
$ws.Columns.Item(7).ColumnWidth = 10
$ws.Columns.Item(15).ColumnWidth = 19.333333333333332
$ws.Columns.Item(16).ColumnWidth = 7.833333333333333
$ws.Columns.Item(17).ColumnWidth = 41.166666666666664
$ws.Columns.Item(25).ColumnWidth = 8.833333333333334

# ---- 3) New "Metadata" worksheet, placed right after VoltageRelay ---------
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "Metadata"

# Helper cell used to coerce plain text (avoids automatic date detection and
# keeps cells free of extra number-format styles) - written then pasted as
# values-only into the real target cell, then cleared again.
$meta.Range("Z1").Formula = "=""Энергохит"""
$meta.Range("Z1").Copy()
$meta.Range("A1").PasteSpecial(-4163)
$meta.Range("Z1").Clear()

$meta.Range("Z1").Formula = "="" 01.08.2012"""
$meta.Range("Z1").Copy()
$meta.Range("C1").PasteSpecial(-4163)
$meta.Range("Z1").Clear()

$meta.Range("Z1").Formula = "=""07.24.2013"""
$meta.Range("Z1").Copy()
$meta.Range("B1").PasteSpecial(-4163)
$meta.Range("Z1").Clear()

$meta.Range("Z1").Formula = "=""Updated prices"""
$meta.Range("Z1").Copy()
$meta.Range("D1").PasteSpecial(-4163)
$meta.Range("Z1").Clear()

# Column widths for the Metadata sheet.
$meta.Columns.Item(1).ColumnWidth = 13
$meta.Columns.Item(2).ColumnWidth = 12.666666666666666
$meta.Columns.Item(3).ColumnWidth = 11.5

# Leave the selection on the Metadata sheet where the author left it.
$meta.Range("E3").Select()

# ---- 4) Restore VoltageRelay as the active sheet/selection -----------------
$ws.Activate()
$ws.Range("P2").Select()
